$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.794.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.20'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.37'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.51%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.27'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.860.33'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.634.59'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0768'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.809.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.09'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.29'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.28%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.124'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.90'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.903'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.135.93'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0157'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.02%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.807'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.769.91'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₆0109'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.25'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.40%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.30%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.25%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.06%  '
